# Apply updated cryptocurrency data (prices / 1h volume % changes)
# Two rows (Bittensor/Cosmos and Kaspa/Dai) were also reordered, which is
# reflected here as their full row contents (name, link, price, volume) swapping.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "545.04",
    # "1.00", "0.999") are not silently coerced into numbers, while
    # resetting back to the Normal style afterwards keeps the cell
    # style index unchanged (matches the un-styled source cells).
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 'D2' '63.621.21'
Set-TextCell 'E2' '  -2.39%  '
Set-TextCell 'D3' '3.330.09'
Set-TextCell 'E3' '  -2.99%  '
Set-TextCell 'E4' '  +0.19%  '
Set-TextCell 'D5' '545.04'
Set-TextCell 'E5' '  -0.78%  '
Set-TextCell 'D6' '171.51'
Set-TextCell 'E6' '  -3.95%  '
Set-TextCell 'E7' '  -4.81%  '
Set-TextCell 'D8' '3.320.95'
Set-TextCell 'E8' '  -2.98%  '
Set-TextCell 'E9' '  +0.01%  '
Set-TextCell 'E10' '  -3.30%  '
Set-TextCell 'E11' '  +1.00%  '
Set-TextCell 'D12' '53.32'
Set-TextCell 'E12' '  -0.43%  '
Set-TextCell 'E13' '  -2.22%  '
Set-TextCell 'E14' '  -3.53%  '
Set-TextCell 'D15' '3.867.22'
Set-TextCell 'E15' '  -2.58%  '
Set-TextCell 'D16' '17.97'
Set-TextCell 'E16' '  -1.93%  '
Set-TextCell 'E17' '  -3.43%  '
Set-TextCell 'D18' '3.312.46'
Set-TextCell 'E18' '  -3.22%  '
Set-TextCell 'E19' '  -1.30%  '
Set-TextCell 'D20' '63.571.35'
Set-TextCell 'E20' '  -2.41%  '
Set-TextCell 'D21' '0.972'
Set-TextCell 'E21' '  -1.13%  '
Set-TextCell 'D22' '409.02'
Set-TextCell 'E22' '  -0.90%  '
Set-TextCell 'E23' '  +0.36%  '
Set-TextCell 'D24' '4.35'
Set-TextCell 'E24' '  +5.73%  '
Set-TextCell 'D25' '13.61'
Set-TextCell 'E25' '  +11.39%  '
Set-TextCell 'E26' '  -3.42%  '
Set-TextCell 'D27' '10.50'
Set-TextCell 'E27' '  -2.02%  '
Set-TextCell 'E28' '  -4.57%  '
Set-TextCell 'D29' '8.55'
Set-TextCell 'E29' '  -4.29%  '
Set-TextCell 'D30' '28.95'
Set-TextCell 'E30' '  -2.65%  '
Set-TextCell 'E31' '  -2.28%  '
Set-TextCell 'B32' 'Bittensor'
Set-TextCell 'C32' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell 'D32' '577.30'
Set-TextCell 'E32' '  -5.02%  '
Set-TextCell 'B33' 'Cosmos'
Set-TextCell 'C33' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 'D33' '11.31'
Set-TextCell 'E33' '  -2.88%  '
Set-TextCell 'E34' '  -2.35%  '
Set-TextCell 'D35' '57.58'
Set-TextCell 'E35' '  -2.35%  '
Set-TextCell 'B36' 'Kaspa'
Set-TextCell 'C36' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell 'D36' '0.147'
Set-TextCell 'E36' '  +1.32%  '
Set-TextCell 'B37' 'Dai'
Set-TextCell 'C37' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 'D37' '1.00'
Set-TextCell 'E37' '  -0.07%  '
Set-TextCell 'D38' '35.05'
Set-TextCell 'E38' '  -5.96%  '
Set-TextCell 'D39' '3.40'
Set-TextCell 'E39' '  +3.50%  '
Set-TextCell 'D40' '0.0₃0734'
Set-TextCell 'E40' '  -6.08%  '
Set-TextCell 'D41' '0.365'
Set-TextCell 'E41' '  -3.45%  '
Set-TextCell 'D42' '3.123.75'
Set-TextCell 'E42' '  -1.94%  '
Set-TextCell 'D43' '0.999'
Set-TextCell 'E43' '  +0.16%  '
Set-TextCell 'E44' '  +0.08%  '
Set-TextCell 'D45' '3.23'
Set-TextCell 'E45' '  +1.63%  '
Set-TextCell 'E46' '  -2.74%  '
Set-TextCell 'D47' '2.41'
Set-TextCell 'E47' '  -4.90%  '
Set-TextCell 'E48' '  -3.92%  '
Set-TextCell 'E49' '  -3.73%  '
Set-TextCell 'D50' '132.44'
Set-TextCell 'E50' '  -3.43%  '
Set-TextCell 'D51' '8.03'
Set-TextCell 'E51' '  -3.91%  '
